$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.580.59"
$ws.Cells.Item(2, 5).Value = "  +7.39%  "
$ws.Cells.Item(3, 4).Value = "1.728.94"
$ws.Cells.Item(3, 5).Value = "  +4.01%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  -0.38%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "332.02"
$ws.Cells.Item(5, 5).Value = "  +0.70%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.9995"
$ws.Cells.Item(6, 5).Value = "  -0.16%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.3735"
$ws.Cells.Item(7, 5).Value = "  +2.35%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3403"
$ws.Cells.Item(8, 5).Value = "  +5.30%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "48.23"
$ws.Cells.Item(9, 5).Value = "  +2.19%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.185"
$ws.Cells.Item(10, 5).Value = "  +4.19%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.07457"
$ws.Cells.Item(11, 5).Value = "  +5.90%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.9989"
$ws.Cells.Item(12, 5).Value = "  -0.38%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "6.424"
$ws.Cells.Item(13, 5).Value = "  +6.09%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "20.12"
$ws.Cells.Item(14, 5).Value = "  +2.99%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.060"
$ws.Cells.Item(15, 5).Value = "  +7.10%  "
$ws.Cells.Item(16, 4).Value = "1.723.32"
$ws.Cells.Item(16, 5).Value = "  +3.75%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001076"
$ws.Cells.Item(17, 5).Value = "  +2.69%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.06653"
$ws.Cells.Item(18, 5).Value = "  +0.90%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "82.28"
$ws.Cells.Item(19, 5).Value = "  +4.82%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.9994"
$ws.Cells.Item(20, 5).Value = "  -0.16%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "16.60"
$ws.Cells.Item(21, 5).Value = "  +5.21%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.208"
$ws.Cells.Item(22, 5).Value = "  +4.86%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "12.76"
$ws.Cells.Item(23, 5).Value = "  +2.58%  "
$ws.Cells.Item(24, 4).Value = "26.544.12"
$ws.Cells.Item(24, 5).Value = "  +7.18%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.444"
$ws.Cells.Item(25, 5).Value = "  -0.23%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.432"
$ws.Cells.Item(26, 5).Value = "  +22.40%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.403"
$ws.Cells.Item(27, 5).Value = "  +0.14%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "150.72"
$ws.Cells.Item(28, 5).Value = "  +1.44%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "19.46"
$ws.Cells.Item(29, 5).Value = "  +4.55%  "
$ws.Cells.Item(30, 4).Value = "1.918.19"
$ws.Cells.Item(30, 5).Value = "  +3.98%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "132.12"
$ws.Cells.Item(31, 5).Value = "  +5.52%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.106"
$ws.Cells.Item(32, 5).Value = "  +0.92%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "6.011"
$ws.Cells.Item(33, 5).Value = "  +6.07%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.08633"
$ws.Cells.Item(34, 5).Value = "  +1.79%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.691"
$ws.Cells.Item(35, 5).Value = "  +3.26%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "12.79"
$ws.Cells.Item(36, 5).Value = "  +6.04%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "5.376"
$ws.Cells.Item(37, 5).Value = "  +4.64%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02344"
$ws.Cells.Item(38, 5).Value = "  +4.04%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.2162"
$ws.Cells.Item(39, 5).Value = "  +4.15%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.06220"
$ws.Cells.Item(40, 5).Value = "  +3.37%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "8.460"
$ws.Cells.Item(41, 5).Value = "  +3.47%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.221"
$ws.Cells.Item(42, 5).Value = "  -0.09%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.6220"
$ws.Cells.Item(43, 5).Value = "  +5.50%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "14.24"
$ws.Cells.Item(44, 5).Value = "  +6.38%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.9987"
$ws.Cells.Item(45, 5).Value = "  -0.22%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.899"
$ws.Cells.Item(46, 5).Value = "  +1.47%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.6029"
$ws.Cells.Item(47, 5).Value = "  +6.96%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "128.78"
$ws.Cells.Item(48, 5).Value = "  +2.94%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.049"
$ws.Cells.Item(49, 5).Value = "  +5.51%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.07185"
$ws.Cells.Item(50, 5).Value = "  +3.22%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "77.17"
$ws.Cells.Item(51, 5).Value = "  +3.45%  "
